$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 43, pushing existing rows 43:46 down to 44:47
$ws.Rows("43:43").Insert()

# Make sure column D on the new row carries the same date number format as the
# rest of the date column.
$ws.Range("D43").NumberFormat = $ws.Range("D44").NumberFormat

# Populate the new row 43 with the new data record.
# Columns A,B,C,E,F,G,H,I,J,K,L,Q,R are constant across this product's rows,
# so copy them from row 44 (the row that used to be row 43 before the insert).
$ws.Range("A43").Value = 10
$ws.Range("B43").Value = "Vega Modelo de Temuco"
$ws.Range("C43").Value = "La Araucanía"
$ws.Range("D43").Value = 44474
$ws.Range("E43").Value = 9
$ws.Range("F43").Value = "Fruta"
$ws.Range("G43").Value = 100108
$ws.Range("H43").Value = "Tropicales y subtropicales"
$ws.Range("I43").Value = 100108007
$ws.Range("J43").Value = "Coco"
$ws.Range("K43").Value = "Sin especificar"
$ws.Range("L43").Value = "Primera"
$ws.Range("M43").Value = 20
$ws.Range("N43").Value = 24000
$ws.Range("O43").Value = 24000
$ws.Range("P43").Value = 24000
$ws.Range("Q43").Value = "$/malla 20 unidades"
$ws.Range("R43").Value = "Perú"
$ws.Range("S43").Value = 1200
$ws.Range("T43").Value = 20
